$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.75"
$ws.Range("E2").Value = "'1.60%"
$ws.Range("D3").Value = "'26.69"
$ws.Range("E3").Value = "'-1.89%"
$ws.Range("D4").Value = "'4.708"
$ws.Range("E4").Value = "'-0.03%"
$ws.Range("D5").Value = "'0.06078"
$ws.Range("E5").Value = "'-1.91%"
$ws.Range("E6").Value = "'0.36%"
$ws.Range("D7").Value = "'0.8499"
$ws.Range("E7").Value = "'-0.10%"
$ws.Range("D8").Value = "'0.9050"
$ws.Range("E8").Value = "'-0.91%"
$ws.Range("D9").Value = "'0.1410"
$ws.Range("E9").Value = "'-0.17%"
$ws.Range("D10").Value = "'0.04866"
$ws.Range("E10").Value = "'7.77%"
$ws.Range("D11").Value = "'0.07088"
$ws.Range("E11").Value = "'0.05%"
$ws.Range("D12").Value = "'0.03178"
$ws.Range("E12").Value = "'1.49%"
$ws.Range("D13").Value = "'0.09022"
$ws.Range("E13").Value = "'-0.22%"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'-0.18%"
$ws.Range("D15").Value = "'0.0006048"
$ws.Range("E15").Value = "'-1.89%"
$ws.Range("D16").Value = "'0.006001"
$ws.Range("E16").Value = "'-0.89%"
$ws.Range("E17").Value = "'-0.10%"
$ws.Range("D18").Value = "'3.168"
$ws.Range("E18").Value = "'0.11%"
$ws.Range("E19").Value = "'3.79%"
$ws.Range("E21").Value = "'-0.79%"
$ws.Range("D22").Value = "'4.077"
$ws.Range("E22").Value = "'-0.05%"
$ws.Range("D23").Value = "'0.04245"
$ws.Range("E23").Value = "'0.32%"
$ws.Range("E24").Value = "'-2.59%"
$ws.Range("D25").Value = "'0.004136"
$ws.Range("E25").Value = "'8.78%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.00%"
$ws.Range("D27").Value = "'0.0001681"
$ws.Range("E27").Value = "'5.01%"
$ws.Range("D40").Value = "'0.03911"
$ws.Range("E40").Value = "'-0.69%"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.03%"
$ws.Range("D42").Value = "'0.004167"
$ws.Range("E42").Value = "'0.94%"
$ws.Range("D43").Value = "'0.002110"
$ws.Range("E43").Value = "'-3.36%"
$ws.Range("E44").Value = "'-8.91%"
$ws.Range("D45").Value = "'0.00005107"
$ws.Range("E45").Value = "'-0.67%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D48").Value = "'0.1563"
$ws.Range("E48").Value = "'-6.20%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.00%"
